$d = $word.ActiveDocument

$replacements = @(
    @{old = "237÷8="; new = "221÷5="},
    @{old = "233÷4="; new = "176÷2="},
    @{old = "978÷9="; new = "490÷9="},
    @{old = "455÷5="; new = "665÷9="},
    @{old = "181÷3="; new = "768÷4="},
    @{old = "128÷6="; new = "728÷6="},
    @{old = "342÷6="; new = "547÷6="},
    @{old = "355÷2="; new = "608÷7="},
    @{old = "279÷3="; new = "937÷7="},
    @{old = "934÷5="; new = "692÷9="},
    @{old = "422÷6="; new = "783÷7="},
    @{old = "914÷2="; new = "754÷4="},
    @{old = "203÷6="; new = "230÷8="},
    @{old = "695÷9="; new = "407÷9="},
    @{old = "304÷2="; new = "186÷2="},
    @{old = "498÷6="; new = "863÷5="},
    @{old = "713÷3="; new = "939÷3="},
    @{old = "447÷9="; new = "892÷8="},
    @{old = "718÷5="; new = "982÷5="},
    @{old = "228÷3="; new = "228÷4="},
    @{old = "360÷8="; new = "209÷2="},
    @{old = "239÷8="; new = "259÷2="},
    @{old = "996÷9="; new = "456÷5="},
    @{old = "912÷6="; new = "532÷4="},
    @{old = "432÷7="; new = "256÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
